$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in previously empty "Union" support cells (column S) with "Y"/"N"
# matching the style already used in the corresponding row's other cells.

# Row 4: S4 -> "Y", keep style 5 (same as it already had)
$ws.Range("S4").Value = "Y"

# Row 5: S5 -> "Y", keep style 5
$ws.Range("S5").Value = "Y"

# Row 6: S6 -> "Y", style becomes 8 (copy format from T6 which already uses style 8)
$ws.Range("S6").Value = "Y"
$ws.Range("T6").Copy()
$ws.Range("S6").PasteSpecial(-4122)

# Row 12: S12 -> "Y", style becomes 8 (copy format from Y12 which already uses style 8)
$ws.Range("S12").Value = "Y"
$ws.Range("Y12").Copy()
$ws.Range("S12").PasteSpecial(-4122)

# Row 13: S13 -> "Y", style becomes 8 (copy format from Y13 which already uses style 8)
$ws.Range("S13").Value = "Y"
$ws.Range("Y13").Copy()
$ws.Range("S13").PasteSpecial(-4122)

# Row 14: J14 changes from "Y" to "N"; S14 -> "N" keeping style 5
$ws.Range("J14").Value = "N"
$ws.Range("S14").Value = "N"

$excel.CutCopyMode = $false

# Update the active selection to H14 as recorded in the saved view state
$ws.Range("H14").Select()
